$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("total")
$ws2 = $wb.Worksheets.Item("sus")

# --- Update data values: column B becomes a per-id sequential timestamp counter ---
# --- (was a flat running row index 1..22) + id 66666 corrected to 666 ---
$ws1.Cells.Item(2,1).Value = 123
$ws1.Cells.Item(2,2).Value = 1
$ws1.Cells.Item(3,1).Value = 679
$ws1.Cells.Item(3,2).Value = 1
$ws1.Cells.Item(4,1).Value = 453
$ws1.Cells.Item(4,2).Value = 1
$ws1.Cells.Item(5,1).Value = 452
$ws1.Cells.Item(5,2).Value = 1
$ws1.Cells.Item(6,1).Value = 123
$ws1.Cells.Item(6,2).Value = 2
$ws1.Cells.Item(7,1).Value = 123
$ws1.Cells.Item(7,2).Value = 3
$ws1.Cells.Item(8,1).Value = 452
$ws1.Cells.Item(8,2).Value = 2
$ws1.Cells.Item(9,1).Value = 123
$ws1.Cells.Item(9,2).Value = 4
$ws1.Cells.Item(10,1).Value = 679
$ws1.Cells.Item(10,2).Value = 2
$ws1.Cells.Item(11,1).Value = 679
$ws1.Cells.Item(11,2).Value = 3
$ws1.Cells.Item(12,1).Value = 123
$ws1.Cells.Item(12,2).Value = 5
$ws1.Cells.Item(13,1).Value = 453
$ws1.Cells.Item(13,2).Value = 2
$ws1.Cells.Item(14,1).Value = 453
$ws1.Cells.Item(14,2).Value = 3
$ws1.Cells.Item(15,1).Value = 452
$ws1.Cells.Item(15,2).Value = 3
$ws1.Cells.Item(16,1).Value = 123
$ws1.Cells.Item(16,2).Value = 6
$ws1.Cells.Item(17,1).Value = 998
$ws1.Cells.Item(17,2).Value = 1
$ws1.Cells.Item(18,1).Value = 453
$ws1.Cells.Item(18,2).Value = 4
$ws1.Cells.Item(19,1).Value = 998
$ws1.Cells.Item(19,2).Value = 2
$ws1.Cells.Item(20,1).Value = 123
$ws1.Cells.Item(20,2).Value = 7
$ws1.Cells.Item(21,1).Value = 998
$ws1.Cells.Item(21,2).Value = 3
$ws1.Cells.Item(22,1).Value = 998
$ws1.Cells.Item(22,2).Value = 4
$ws1.Cells.Item(23,1).Value = 666
$ws1.Cells.Item(23,2).Value = 1

# --- Append 2 new observations for id 343 ---
$ws1.Cells.Item(24,1).Value = 343
$ws1.Cells.Item(24,2).Value = 1
$ws1.Cells.Item(24,3).Value = 63.2
$ws1.Cells.Item(24,4).Value = 71
$ws1.Cells.Item(24,5).Value = 29
$ws1.Cells.Item(25,1).Value = 343
$ws1.Cells.Item(25,2).Value = 2
$ws1.Cells.Item(25,3).Value = 50
$ws1.Cells.Item(25,4).Value = 77
$ws1.Cells.Item(25,5).Value = 18

# --- Colour-code column A by id, reusing one theme fill colour per distinct id ---
$ws1.Cells.Item(2,1).Interior.ThemeColor = 10
$ws1.Cells.Item(3,1).Interior.ThemeColor = 9
$ws1.Cells.Item(4,1).Interior.ThemeColor = 8
$ws1.Cells.Item(5,1).Interior.ThemeColor = 6
$ws1.Cells.Item(6,1).Interior.ThemeColor = 10
$ws1.Cells.Item(7,1).Interior.ThemeColor = 10
$ws1.Cells.Item(8,1).Interior.ThemeColor = 6
$ws1.Cells.Item(9,1).Interior.ThemeColor = 10
$ws1.Cells.Item(10,1).Interior.ThemeColor = 9
$ws1.Cells.Item(11,1).Interior.ThemeColor = 9
$ws1.Cells.Item(12,1).Interior.ThemeColor = 10
$ws1.Cells.Item(13,1).Interior.ThemeColor = 8
$ws1.Cells.Item(14,1).Interior.ThemeColor = 8
$ws1.Cells.Item(15,1).Interior.ThemeColor = 6
$ws1.Cells.Item(16,1).Interior.ThemeColor = 10
$ws1.Cells.Item(17,1).Interior.ThemeColor = 4
$ws1.Cells.Item(18,1).Interior.ThemeColor = 8
$ws1.Cells.Item(19,1).Interior.ThemeColor = 4
$ws1.Cells.Item(20,1).Interior.ThemeColor = 10
$ws1.Cells.Item(21,1).Interior.ThemeColor = 4
$ws1.Cells.Item(22,1).Interior.ThemeColor = 4
$ws1.Cells.Item(23,1).Interior.ThemeColor = 2
$ws1.Cells.Item(24,1).Interior.ThemeColor = 4
$ws1.Cells.Item(25,1).Interior.ThemeColor = 4

# --- Restore each sheet's selection / cursor position ---
$ws2.Range("A32").Select()
$ws1.Range("C25:E25").Select()

Write-Host "edit complete"